# Fixes a one-row misalignment in the Brazilian-states table:
#  - Row 2 ("TO" row) actually held the Brasil grand-total figures; it is
#    corrected to the "*Tot"/"00"/"Brasil" total row with the real totals.
#  - Rows 9-28 were each showing the state-name/figures that belonged to the
#    row below them (the "Tocantins" state was missing entirely); every row
#    from 9 to 28 is shifted to the correct coduf/uf/figures, and a new
#    "Tocantins" (TO, coduf 17) row appears at row 9.
#  - Row 29 (previously #N/A/#N/A) becomes the proper "DF"/"53"/"Distrito
#    Federal" row.
# Columns: A=sguf, B=coduf (text, keeps leading zeros), C=uf (full name),
#          D=Num. Dep. Fed, E=Eleitos 14 e perc, F=Não eleitos e perc

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  A = "*Tot"; B = "00"; C = "Brasil";               D = 513; E = 468; F = 45 },
    @{ Row = 9;  A = "TO";   B = "17"; C = "Tocantins";             D = 8;   E = 7;   F = 1  },
    @{ Row = 10; A = "MA";   B = "21"; C = "Maranhão";              D = 18;  E = 14;  F = 4  },
    @{ Row = 11; A = "PI";   B = "22"; C = "Piauí";                 D = 10;  E = 9;   F = 1  },
    @{ Row = 12; A = "CE";   B = "23"; C = "Ceará";                 D = 22;  E = 22;  F = 0  },
    @{ Row = 13; A = "RN";   B = "24"; C = "Rio Grande do Norte";   D = 8;   E = 8;   F = 0  },
    @{ Row = 14; A = "PB";   B = "25"; C = "Paraíba";               D = 12;  E = 12;  F = 0  },
    @{ Row = 15; A = "PE";   B = "26"; C = "Pernambuco";            D = 25;  E = 24;  F = 1  },
    @{ Row = 16; A = "AL";   B = "27"; C = "Alagoas";               D = 9;   E = 8;   F = 1  },
    @{ Row = 17; A = "SE";   B = "28"; C = "Sergipe";               D = 8;   E = 8;   F = 0  },
    @{ Row = 18; A = "BA";   B = "29"; C = "Bahia";                 D = 39;  E = 38;  F = 1  },
    @{ Row = 19; A = "MG";   B = "31"; C = "Minas Gerais";          D = 53;  E = 48;  F = 5  },
    @{ Row = 20; A = "ES";   B = "32"; C = "Espírito Santo";        D = 10;  E = 8;   F = 2  },
    @{ Row = 21; A = "RJ";   B = "33"; C = "Rio de Janeiro";        D = 46;  E = 43;  F = 3  },
    @{ Row = 22; A = "SP";   B = "35"; C = "São Paulo";             D = 70;  E = 62;  F = 8  },
    @{ Row = 23; A = "PR";   B = "41"; C = "Paraná";                D = 30;  E = 27;  F = 3  },
    @{ Row = 24; A = "SC";   B = "42"; C = "Santa Catarina";        D = 16;  E = 14;  F = 2  },
    @{ Row = 25; A = "RS";   B = "43"; C = "Rio Grande do Sul";     D = 31;  E = 28;  F = 3  },
    @{ Row = 26; A = "MS";   B = "50"; C = "Mato Grosso do Sul";    D = 8;   E = 7;   F = 1  },
    @{ Row = 27; A = "MT";   B = "51"; C = "Mato Grosso";           D = 8;   E = 7;   F = 1  },
    @{ Row = 28; A = "GO";   B = "52"; C = "Goiás";                 D = 17;  E = 16;  F = 1  },
    @{ Row = 29; A = "DF";   B = "53"; C = "Distrito Federal";      D = 8;   E = 6;   F = 2  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
